# Update betting odds values for row 5 and row 7 on Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 5 changes
$ws.Range("Q5").Value = 2.4
$ws.Range("R5").Value = 1.53

# Row 7 changes
$ws.Range("G7").Value = 2.5
$ws.Range("I7").Value = 2.75
$ws.Range("J7").Value = 3
$ws.Range("K7").Value = 2.05
$ws.Range("L7").Value = 3.35
$ws.Range("S7").Value = 1.4
$ws.Range("T7").Value = 2.52
$ws.Range("W7").Value = 7.7
$ws.Range("Y7").Value = 9.5
$ws.Range("AA7").Value = 21
$ws.Range("AB7").Value = 32
$ws.Range("AH7").Value = 7.7
$ws.Range("AI7").Value = 13
$ws.Range("AJ7").Value = 10.25
$ws.Range("AL7").Value = 26
$ws.Range("AM7").Value = 37
$ws.Range("AN7").Value = 4.35
$ws.Range("AO7").Value = 13
$ws.Range("AP7").Value = 20
$ws.Range("AQ7").Value = 55
$ws.Range("AR7").Value = 80
$ws.Range("AS7").Value = 250
$ws.Range("AT7").Value = 2.5
$ws.Range("AU7").Value = 6.8
$ws.Range("AW7").Value = 4.65
$ws.Range("AX7").Value = 15
$ws.Range("AY7").Value = 23
$ws.Range("AZ7").Value = 70
$ws.Range("BA7").Value = 110

$wb.Save()
